$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the two new columns for "age" (E) and "gender" (F) right after
#    "email" (D). Everything from the old "address" column onward shifts
#    right by two (old E:L -> new G:N).
# ---------------------------------------------------------------------------
$ws.Columns("E:F").Insert()

# ---------------------------------------------------------------------------
# 2. Insert the new "event_id" column before "date_of_initiation".
#    After the first insert, date_of_initiation lives in column N, so we
#    insert a fresh column there, pushing it to O.
# ---------------------------------------------------------------------------
$ws.Columns("N:N").Insert()

# ---------------------------------------------------------------------------
# 3. Header row
# ---------------------------------------------------------------------------
$ws.Range("E1").Value = "age"
$ws.Range("F1").Value = "gender"
$ws.Range("N1").Value = "event_id"

# Header style for the new event_id header: bold + bordered (like the other
# headers) but centered vertically instead of top-aligned.
$ws.Range("N1").Font.Bold = $true
$ws.Range("N1").Borders.LineStyle = 1
$ws.Range("N1").HorizontalAlignment = -4108
$ws.Range("N1").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 4. Data rows (2-11): fill age / gender / event_id, and overwrite the old
#    state/district/country text columns (now H/I/J) with the new numeric
#    ids used by the upload template.
# ---------------------------------------------------------------------------
$ages    = @(23, 32, 27, 33, 34, 19, 24, 23, 25, 26)
$genders = @("Male", "Male", "Male", "Female", "Male", "Female", "Male", "Female", "Male", "Female")
$eventId = @(1, 1, 1, 2, 2, 2, 2, 2, 2, 2)

for ($i = 0; $i -lt 10; $i++) {
    $r = $i + 2

    $ws.Range("E$r").Value = $ages[$i]
    $ws.Range("F$r").Value = $genders[$i]

    $ws.Range("H$r").Value = 1682
    $ws.Range("I$r").Value = 48723
    $ws.Range("J$r").Value = 101

    $ws.Range("N$r").Value = $eventId[$i]
}

# ---------------------------------------------------------------------------
# 5. Styles: center/center alignment on id, age, gender, event_id columns;
#    center-only alignment on instructor_id.
# ---------------------------------------------------------------------------
$ws.Range("A2:A11").HorizontalAlignment = -4108
$ws.Range("A2:A11").VerticalAlignment = -4108

$ws.Range("E2:F11").HorizontalAlignment = -4108
$ws.Range("E2:F11").VerticalAlignment = -4108

$ws.Range("M2:M11").HorizontalAlignment = -4108

$ws.Range("N2:N11").HorizontalAlignment = -4108
$ws.Range("N2:N11").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 6. Column widths
# ---------------------------------------------------------------------------
$ws.Columns("A:C").ColumnWidth = 14.88
$ws.Columns("D:D").ColumnWidth = 19.74
$ws.Columns("E:K").ColumnWidth = 14.88
$ws.Columns("L:L").ColumnWidth = 27.31
$ws.Columns("M:N").ColumnWidth = 14.88
$ws.Columns("O:O").ColumnWidth = 16.45
$ws.Columns("P:P").ColumnWidth = 14.88

# ---------------------------------------------------------------------------
# 7. Misc sheet cosmetics to mirror the authored workbook.
# ---------------------------------------------------------------------------
$ws.Range("D17").Select()

Write-Output "edit applied"
